$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the crypto-price refresh diff.
# D-column price values are forced to Text (leading apostrophe, Excel's
# standard "treat as text" convention) and the style is reset to "Normal"
# afterwards so numeric-looking strings (e.g. "213.26", "1.70", "0.0486")
# keep their exact original formatting instead of being re-parsed as numbers.

$ws.Range('D2').Value = "'27.719.42"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.20%  '

$ws.Range('D3').Value = "'1.645.98"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.62%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').Value = "'213.26"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.03%  '

$ws.Range('E6').Value = '  +3.69%  '

$ws.Range('E7').Value = '  +0.06%  '

$ws.Range('D8').Value = "'23.11"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.04%  '

$ws.Range('D9').Value = "'0.259"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.22%  '

$ws.Range('E10').Value = '  +0.01%  '

$ws.Range('E11').Value = '  +1.84%  '

$ws.Range('D12').Value = "'1.878.88"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.62%  '

$ws.Range('D13').Value = "'1.643.59"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.70%  '

$ws.Range('E14').Value = '  -0.80%  '

$ws.Range('D15').Value = "'0.562"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.11%  '

$ws.Range('D16').Value = "'64.29"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.03%  '

$ws.Range('D17').Value = "'27.694.60"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.10%  '

$ws.Range('D18').Value = "'231.00"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.20%  '

$ws.Range('E19').Value = '  +0.09%  '

$ws.Range('E20').Value = '  +2.61%  '

$ws.Range('E21').Value = '  +0.03%  '

$ws.Range('E22').Value = '  -1.05%  '

$ws.Range('D23').Value = "'10.03"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +7.32%  '

$ws.Range('E24').Value = '  -3.07%  '

$ws.Range('D25').Value = "'149.58"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.37%  '

$ws.Range('E26').Value = '  -1.84%  '

$ws.Range('E27').Value = '  +1.23%  '

$ws.Range('D28').Value = "'15.68"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.33%  '

$ws.Range('E29').Value = '  +0.06%  '

$ws.Range('E30').Value = '  -0.12%  '

$ws.Range('D31').Value = "'0.0486"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.13%  '

$ws.Range('E32').Value = '  +0.16%  '

$ws.Range('E33').Value = '  +1.23%  '

$ws.Range('D34').Value = "'1.443.09"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.75%  '

$ws.Range('E35').Value = '  +2.01%  '

$ws.Range('E36').Value = '  -1.12%  '

$ws.Range('D37').Value = "'0.570"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.31%  '

$ws.Range('E38').Value = '  -2.52%  '

$ws.Range('E39').Value = '  -0.87%  '

$ws.Range('D40').Value = "'0.901"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +13.97%  '

$ws.Range('D41').Value = "'1.03"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.42%  '

$ws.Range('E42').Value = '  +0.07%  '

$ws.Range('D43').Value = "'5.66"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.61%  '

$ws.Range('E44').Value = '  -0.36%  '

$ws.Range('E45').Value = '  +1.92%  '

$ws.Range('D46').Value = "'65.62"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.05%  '

$ws.Range('D47').Value = "'1.788.17"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.54%  '

$ws.Range('D48').Value = "'1.70"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.37%  '

$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = "'0.0₆0108"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.52%  '

$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').Value = "'85.94"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.36%  '

$ws.Range('D51').Value = "'0.0991"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.92%  '
